# Daily scrape update - 2026-01-06
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width adjustments (C, D, H got narrower) ---
# Excel's ColumnWidth setter stores width + 5/6 char in the OOXML <col width>,
# so back the requested value off by 5/6 to land on the exact target.
$ws.Columns.Item(3).ColumnWidth = (49 - 5/6)
$ws.Columns.Item(4).ColumnWidth = (35 - 5/6)
$ws.Columns.Item(8).ColumnWidth = (45 - 5/6)

# --- Row data (rows 2-9), columns A-H ---
$data = @(
    @("1330365","https://aiesec.org/opportunity/global-talent/1330365","Creative Styling & Brand Experience Intern","Hyderabad, Telangana, India","No","1 applicant","9 - 12 Weeks","MPF clothing collection PVT LTD"),
    @("1330185","https://aiesec.org/opportunity/global-talent/1330185","Front-End Web Developer","Ciudad Juárez, Chihuahua, Mexico","No","23 applicants","6 - 18 Months","EP&O Corporation"),
    @("1329279","https://aiesec.org/opportunity/global-talent/1329279","Markets Commercial Ops trainee","Bruxelles, Belgio","No","101 applicants","6 - 18 Months","UCB"),
    @("1328614","https://aiesec.org/opportunity/global-talent/1328614","Field Service Engineer [EU Preferred]","Madrid, Spain","No","140 applicants","6 - 18 Months","Mitsubishi Power Europe Sucursal en España"),
    @("1326310","https://aiesec.org/opportunity/global-talent/1326310","Back - End Developer","Glyfada, Greece","No","284 applicants","3 - 6 Months","Validata Software"),
    @("1325118","https://aiesec.org/opportunity/global-talent/1325118","IT Support Specialist (Flexible RE dates)","Χολαργός, Ελλάδα","No","151 applicants","6 - 18 Months","WizzIT"),
    @("1324106","https://aiesec.org/opportunity/global-talent/1324106","Tourism Specialist - Intern","Nugegoda, Sri Lanka","No","23 applicants","9 - 12 Weeks","Brand Corridor (Pvt) Ltd"),
    @("1319024","https://aiesec.org/opportunity/global-talent/1319024","[Impact Brazil] - A. I. Technologies Developer","Uberlândia, MG, Brasil","No","152 applicants","6 - 18 Months","Neospace A. I. Technologies")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $vals = $data[$i]

    # Column A holds an opportunity id that looks numeric; force it to stay
    # text (matches the source file's inlineStr) and strip the quote-prefix
    # style Excel applies so the cell keeps the workbook's default style.
    $cellA = $ws.Cells.Item($row, 1)
    $cellA.NumberFormat = "@"
    $cellA.Value = $vals[0]
    $cellA.Style = $ws.Cells.Item(1, 5).Style

    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
    $ws.Cells.Item($row, 6).Value = $vals[5]
    $ws.Cells.Item($row, 7).Value = $vals[6]
    $ws.Cells.Item($row, 8).Value = $vals[7]
}
